$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append new mail-log row (A11:G11) ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A11").Value = "Is product X op voorraad?"
$logs.Range("B11").Value = "mailmind.test@zohomail.eu"
$logs.Range("C11").Value = "Ik ben geïnteresseerd in product X. Is dit momenteel op voorraad?"
$logs.Range("D11").Value = "Productinformatie"
$logs.Range("F11").Value = "2025-06-19 21:20:13"
$logs.Range("G11").Value = "Nee"

# Extend the conditional formatting ranges so they keep covering the new row
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D11"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G11"))

# --- "Dashboard" sheet: swap the top two category rows (counts updated) ---
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A2").Value = "Productinformatie"
$dashboard.Range("B2").Value = 3
$dashboard.Range("A3").Value = "Samenwerking / Partnerverzoek"
$dashboard.Range("B3").Value = 2
